$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.934.45"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.637.92"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +1.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.64"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.865.50"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "1.636.04"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "0.0₃0758"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.54"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "25.962.74"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.63"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.90"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.88"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "1.138.79"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.45"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.800"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("E44").Value = "  -4.30%  "
$ws.Range("D45").Value = "1.774.43"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +5.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.59"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.415"
$ws.Range("E51").Value = "  -0.27%  "
